$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Price"
$ws.Range("D1").Value = "Quantity"
